$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Test Case Title column (F) values exactly as per corrected test data
$ws.Range("F2").Value = "loginTest"
$ws.Range("F3").Value = "loginTest2"
$ws.Range("F4").Value = "Homepage1"
$ws.Range("F5").Value = "Homepage2"
$ws.Range("F6").Value = "Homepage3"
$ws.Range("F7").Value = "Homepage4"
$ws.Range("F8").Value = "Homepage5"

# Update the active cell selection on the sheet
$ws.Range("F9").Select()
